# Re-grade the "CustomerMappingTest Class" / "100% passing of all the test
# cases" row: the grader found the submission failed because of an
# incorrect existence-check / missing hashmap insert (not a missing test
# file), and awarded partial credit instead of zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Points awarded for row 34 (was 0, now 5 out of 7)
$ws.Range("E34").Value = 5

# Updated grading comment for row 34
$ws.Range("F34").Value = "(-5) For incorrect condition for checking to check customer exists or not and not adding new customer into hash map."

# Reflect where the grader's cursor ended up after the edit
$ws.Range("E41").Select()
